$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")

# Header for new column G
$ws.Range("G1").Value = "Invoice Number"
$ws.Range("G1").Font.Size = 14
$ws.Rows.Item(1).RowHeight = 19

# Invoice numbers F2F456 .. F2F473 for rows 2..19
$invoiceNumbers = @(
    "F2F456", "F2F457", "F2F458", "F2F459", "F2F460",
    "F2F461", "F2F462", "F2F463", "F2F464", "F2F465",
    "F2F466", "F2F467", "F2F468", "F2F469", "F2F470",
    "F2F471", "F2F472", "F2F473"
)

for ($i = 0; $i -lt $invoiceNumbers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $invoiceNumbers[$i]
}

# Adjust column G width (closest achievable value given pixel-quantized COM rounding)
$ws.Columns.Item(7).ColumnWidth = 15.33

# Update selection to match the recorded edit
[void]$ws.Range("F11").Select()
